$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "[Michael%Irvine%NULL%0, Daniel%Coombs%NULL%1, Julianne%Skarha%NULL%1, Brandon%del Pozo%NULL%1, Josiah%Rich%NULL%1, Faye%Taxman%NULL%1, Traci C.%Green%Traci.c.green@gmail.com%1]"
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = "Springer US"
$ws.Range("E3").Value = "[Shaun%Truelove%NULL%0, Orit%Abrahim%NULL%2, Orit%Abrahim%NULL%0, Chiara%Altare%NULL%2, Chiara%Altare%NULL%0, Stephen A.%Lauer%NULL%2, Stephen A.%Lauer%NULL%0, Krya H.%Grantz%NULL%2, Krya H.%Grantz%NULL%0, Andrew S.%Azman%NULL%2, Andrew S.%Azman%NULL%0, Paul%Spiegel%NULL%2, Paul%Spiegel%NULL%0, Parveen%Parmar%NULL%3, Parveen%Parmar%NULL%0, Parveen%Parmar%NULL%0, NULL%NULL%NULL%0, NULL%NULL%NULL%0, NULL%NULL%NULL%0]"
$ws.Range("I3").Value = ""
$ws.Range("J3").Value = "Public Library of Science"
$ws.Range("C4").Value = "Unknown Title"
$ws.Range("D4").Value = "Unknown Abstract"
$ws.Range("E4").Value = "[]"
$ws.Range("F4").Value = "not found"
$ws.Range("G4").Value = "N/A"
$ws.Range("H4").Value = "1970-01-01"
$ws.Range("I4").Value = ""
$ws.Range("C5").Value = "Unknown Title"
$ws.Range("D5").Value = "Unknown Abstract"
$ws.Range("E5").Value = "[]"
$ws.Range("F5").Value = "not found"
$ws.Range("G5").Value = "N/A"
$ws.Range("H5").Value = "1970-01-01"
$ws.Range("I5").Value = ""
$ws.Range("D6").Value = "The economic and health consequences of the COVID-19 pandemic pose a particular threat to vulnerable groups, such as migrants, particularly forcibly displaced populations.
 The aim of this review is (i) to synthesize the evidence on risk of infection and transmission among migrants, refugees, asylum seekers and internally displaced populations, and (ii) the effect of lockdown measures on these populations.
 We searched MEDLINE and WOS, preprint servers, and pertinent websites between 1st December 2019 and 26th June 2020. The included studies showed a high heterogeneity in study design, population, outcome and quality.
 The incidence risk of SARS-CoV-2 varied from 0•12% to 2•08% in non-outbreak settings and from 5•64% to 21•15% in outbreak settings.
 Migrants showed a lower hospitalization rate compared to non-migrants.
 Negative impacts on mental health due to lockdown measures were found across respective studies.
 However, findings show a tenuous and heterogeneous data situation, showing the need for more robust and comparative study designs.
"
$ws.Range("E6").Value = "[Maren%Hintermeier%NULL%0, Hande%Gencer%NULL%1, Katja%Kajikhina%NULL%1, Sven%Rohleder%NULL%1, Claudia%Hövener%NULL%1, Marie%Tallarek%NULL%1, Jacob%Spallek%NULL%1, Kayvan%Bozorgmehr%kayvan.bozorgmehr@uni-bielefeld.de%1]"
$ws.Range("I6").Value = ""
$ws.Range("J6").Value = "Elsevier"
$ws.Range("E7").Value = "[David%Koh%NULL%0]"
$ws.Range("I7").Value = ""
$ws.Range("J7").Value = "BMJ Publishing Group"
$ws.Range("E8").Value = "[MH%Chew%chew.min.hoe@singhealth.com.sg%0, F.H.%Koh%NULL%1, JT%Wu%NULL%1, S.%Ngaserin%NULL%1, A.%Ng%NULL%1, BC%Ong%NULL%1, V.J.%Lee%NULL%1]"
$ws.Range("I8").Value = ""
$ws.Range("J8").Value = "The Healthcare Infection Society. Published by Elsevier Ltd."
Write-Host "Applied fold_4/98.xlsx citation graph fixes"
